$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.133.23'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.656.69'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.49%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5247'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.48%  '
$ws.Range("E7").Value = '  -0.46%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2620'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.75%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06291'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.55'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07808'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.486'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.668.37'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.885.08'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5554'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.38%  '
$ws.Range("D16").Value = '0.0₅8023'
$ws.Range("E16").Value = '  -2.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.99'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.145.58'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("E19").Value = '  -0.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.630'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '195.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.08'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.958'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.33%  '
$ws.Range("E24").Value = '  -0.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.07'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1206'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.170'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.39%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.94'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.492'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05695'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.269'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.486'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.343'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.587'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.59%  '
$ws.Range("E35").Value = '  -0.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9503'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.416'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5706'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01596'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.939'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.064.48'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8454'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.12%  '
$ws.Range("E43").Value = '  -0.43%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '103.46'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.795.83'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.72'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.006'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.56%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4401'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05304'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.79%  '
$ws.Range("E50").Value = '  -3.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.981'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.34%  '
